# Update cryptocurrency price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain decimal number (e.g. "520.17")
# would otherwise be auto-converted from text to a numeric value by Excel.
# Pre-formatting them as Text preserves the original inline-string/text type.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '57.058.08'
$ws.Range("E2").Value = '  -2.25%  '
$ws.Range("D3").Value = '3.062.95'
$ws.Range("E3").Value = '  -2.51%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '520.17'
$ws.Range("E5").Value = '  -2.51%  '
$ws.Range("D6").Value = '135.19'
$ws.Range("E6").Value = '  -5.18%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.060.52'
$ws.Range("E8").Value = '  -2.55%  '
$ws.Range("D9").Value = '0.470'
$ws.Range("E9").Value = '  +5.63%  '
$ws.Range("D10").Value = '7.25'
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("D11").Value = '0.106'
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("E12").Value = '  +2.09%  '
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("D14").Value = '3.583.66'
$ws.Range("E14").Value = '  -2.69%  '
$ws.Range("D15").Value = '25.02'
$ws.Range("E15").Value = '  -2.45%  '
$ws.Range("E16").Value = '  -4.07%  '
$ws.Range("D17").Value = '57.090.57'
$ws.Range("E17").Value = '  -2.22%  '
$ws.Range("D18").Value = '3.067.17'
$ws.Range("E18").Value = '  -2.25%  '
$ws.Range("D20").Value = '12.37'
$ws.Range("E20").Value = '  -3.50%  '
$ws.Range("D21").Value = '7.75'
$ws.Range("E21").Value = '  -2.91%  '
$ws.Range("D22").Value = '347.64'
$ws.Range("E22").Value = '  +1.20%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '68.70'
$ws.Range("E24").Value = '  +1.49%  '
$ws.Range("D25").Value = '0.496'
$ws.Range("E25").Value = '  -3.36%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = '0.164'
$ws.Range("E27").Value = '  -3.26%  '
$ws.Range("D28").Value = '0.0₃0841'
$ws.Range("E28").Value = '  -10.20%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '7.11'
$ws.Range("E30").Value = '  -4.22%  '
$ws.Range("E31").Value = '  -2.57%  '
$ws.Range("D32").Value = '20.83'
$ws.Range("E32").Value = '  -1.34%  '
$ws.Range("D33").Value = '5.73'
$ws.Range("E33").Value = '  -10.77%  '
$ws.Range("E34").Value = '  -0.35%  '
$ws.Range("D35").Value = '157.45'
$ws.Range("E35").Value = '  -0.36%  '
$ws.Range("E36").Value = '  -6.52%  '
$ws.Range("D37").Value = '5.95'
$ws.Range("E37").Value = '  -4.52%  '
$ws.Range("D38").Value = '25.11'
$ws.Range("E38").Value = '  -4.65%  '
$ws.Range("E39").Value = '  -2.89%  '
$ws.Range("E40").Value = '  -2.58%  '
$ws.Range("D41").Value = '1.56'
$ws.Range("E41").Value = '  -6.44%  '
$ws.Range("D42").Value = '4.00'
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("E43").Value = '  -3.13%  '
$ws.Range("D44").Value = '2.404.97'
$ws.Range("E44").Value = '  +5.53%  '
$ws.Range("D45").Value = '36.43'
$ws.Range("E45").Value = '  -0.75%  '
$ws.Range("D47").Value = '3.098.49'
$ws.Range("E47").Value = '  -2.65%  '
$ws.Range("E48").Value = '  -1.75%  '
$ws.Range("D49").Value = '5.95'
$ws.Range("E49").Value = '  -2.36%  '
$ws.Range("D50").Value = '0.923'
$ws.Range("E50").Value = '  -8.78%  '
$ws.Range("D51").Value = '19.28'
$ws.Range("E51").Value = '  -6.80%  '
